$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A20:E20").NumberFormat = "@"
$ws.Range("A20").Value = "2024-06-13"
$ws.Range("B20").Value = "신한글로벌액티브리츠"
$ws.Range("C20").Value = "신한, 한국"
$ws.Range("D20").Value = "2024-06-18"
$ws.Range("E20").Value = "2024-07-01"
$ws.Range("F20").Value = 70000002
$ws.Range("G20").Value = 23333334
$ws.Range("H20").Value = "-"
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 3800
$ws.Range("K20").Value = "-"
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = "-"
$ws.Range("N20").Value = "-"
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = "-"
$ws.Range("Q20").Value = "-"
$ws.Range("R20").Value = "120.3 : 1"
$ws.Range("S20").Value = "-"
$ws.Range("T20").Value = "-"

$ws.Range("A20:E20").Style = "Normal"
